$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'хлеб',
    'вино',
    'скот',
    'холст',
    'кожа',
    'мед',
    'пиво',
    'сукно',
    'овчина',
    'лошадь',
    'воск',
    'масло',
    'сало',
    'железо',
    'Крымскую соль',
    'колеса',
    'полотно',
    'сено',
    'говядина',
    'парча',
    'позумент',
    'табак',
    'шелк',
    'сахар',
    'выбойка',
    'чулок',
    'лес',
    'лыко',
    'ладан',
    'сани',
    'китайка',
    'сапог',
    'коса',
    'горшок',
    'платок',
    'ром',
    'овца',
    'конь',
    'рогожа',
    'гвоздь',
    'замок',
    'обод',
    'веревка',
    'сосуд',
    'дуга',
    'брусья',
    'скотский кожа',
    'хомут',
    'роза',
    'гумми',
    'покроми',
    'котел',
    'нитка',
    'сковорода',
    'бечева'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
